$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 207 entirely; this shifts all subsequent rows up by one.
$ws.Rows.Item(207).Delete()
